$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 153, shifting existing rows 153:251 down to 154:252
$ws.Rows("153").Insert()

# Populate the newly inserted row with the new daily record
$ws.Range("A153").Value = 8
$ws.Range("B153").Value = "Terminal La Palmera de La Serena"
$ws.Range("C153").Value = "Coquimbo"
$ws.Range("D153").Value = 44582
$ws.Range("E153").Value = 4
$ws.Range("F153").Value = 100112032
$ws.Range("G153").Value = "Zapallo italiano"
$ws.Range("H153").Value = "Sin especificar"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 500
$ws.Range("K153").Value = 10000
$ws.Range("L153").Value = 11000
$ws.Range("M153").Value = 10500
$ws.Range("N153").Value = "$/caja 70 unidades"
$ws.Range("O153").Value = "Provincia de Limarí"
$ws.Range("P153").Value = 150
$ws.Range("Q153").Value = 70
$ws.Range("R153").Value = "Hortaliza"
